$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28, shifting rows 28-30 down to 29-31
$ws.Rows.Item(28).Insert()

# Update row 27 (exposure period change)
$ws.Range("C27").Value = "28/12/2020 10:30pm-12.00am"

# Fill in the new row 28 with the Nandos entry
$ws.Range("A28").Value = "Melbourne"
$ws.Range("B28").Value = "Nandos  27 Elizabeth Street, Melbourne"
$ws.Range("C28").Value = "01/01/2021 1:00am - 2:00am"
$ws.Range("D28").Value = "Case dined at venue"
